$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ADAM SYAFIQ BIN SAFIAN"
$ws.Range("C2").Value = "011-25686811"
$ws.Range("D2").Value = "adamsyafiq31@gmail.com"

$ws.Range("G2").Value = "[{'job_title': 'Document Controller', 'job_company': 'China Communication Construction (ECRL) Sdn Bhd', 'Industries': 'N/A', 'start_date': '2022-01', 'end_date': '2024-05-21 02:50:37.826879', 'job_location': 'N/A', 'job_duration': '2 years 4 months'}]"

$ws.Range("I2").Value = "[{'field_of_study': 'Bachelor in Engineering Technology (Facilities Maintenance Engineering)', 'level': ""Bachelor's Degree"", 'cgpa': 'N/A', 'university': 'UniKL Mitec', 'start_date': '2018', 'year_of_graduation': '2021'}, {'field_of_study': 'Diploma in Construction Technology (Building Services and Maintenance)', 'level': 'Diploma', 'cgpa': 'N/A', 'university': 'KKTM Sri Gading', 'start_date': '2015', 'year_of_graduation': '2018'}]"

$ws.Range("J2").Value = "['N/A']"

$ws.Range("K2").Value = "['Teamwork', 'Time Management', 'Leadership', 'Microsoft Office', 'Bahasa Malaysia']"

$ws.Range("L2").Value = "['Mandarin', 'English']"
